$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Fgf2"
$ws.Cells.Item(2,3).Value = "Gpc4"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 0.747119
$ws.Cells.Item(2,8).Value = 2.241357
$ws.Cells.Item(2,9).Value = 0.03096954854571248
$ws.Cells.Item(2,10).Value = 0.03096954854571248
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 6.992380999999999
$ws.Cells.Item(2,14).Value = 20.977143
$ws.Cells.Item(2,15).Value = 0.1454502834669897
$ws.Cells.Item(2,16).Value = 0.1454502834669897
$ws.Cells.Item(2,17).Value = 5.224140700338999
$ws.Cells.Item(2,18).Value = 47.017266303051
$ws.Cells.Item(2,19).Value = 0.004504529614818578
$ws.Cells.Item(2,20).Value = 0.004504529614818578

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Fgf2"
$ws.Cells.Item(3,3).Value = "Gpc4"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 0.747119
$ws.Cells.Item(3,8).Value = 2.241357
$ws.Cells.Item(3,9).Value = 0.03096954854571248
$ws.Cells.Item(3,10).Value = 0.03096954854571248
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 30.83466
$ws.Cells.Item(3,14).Value = 92.50398
$ws.Cells.Item(3,15).Value = 0.641399551541635
$ws.Cells.Item(3,16).Value = 0.641399551541635
$ws.Cells.Item(3,17).Value = 23.03716034454
$ws.Cells.Item(3,18).Value = 207.33444310086
$ws.Cells.Item(3,19).Value = 0.01986385454866688
$ws.Cells.Item(3,20).Value = 0.01986385454866688

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Fgf2"
$ws.Cells.Item(4,3).Value = "Gpc4"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 0.747119
$ws.Cells.Item(4,8).Value = 2.241357
$ws.Cells.Item(4,9).Value = 0.03096954854571248
$ws.Cells.Item(4,10).Value = 0.03096954854571248
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 10.246987
$ws.Cells.Item(4,14).Value = 30.740961
$ws.Cells.Item(4,15).Value = 0.2131501649913754
$ws.Cells.Item(4,16).Value = 0.2131501649913754
$ws.Cells.Item(4,17).Value = 7.655718680452999
$ws.Cells.Item(4,18).Value = 68.90146812407698
$ws.Cells.Item(4,19).Value = 0.006601164382227024
$ws.Cells.Item(4,20).Value = 0.006601164382227024

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Fgf2"
$ws.Cells.Item(5,3).Value = "Gpc4"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 19.74619233333334
$ws.Cells.Item(5,8).Value = 59.23857700000001
$ws.Cells.Item(5,9).Value = 0.8185184181638298
$ws.Cells.Item(5,10).Value = 0.8185184181638298
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 0.6666666666666666
$ws.Cells.Item(5,13).Value = 6.992380999999999
$ws.Cells.Item(5,14).Value = 20.977143
$ws.Cells.Item(5,15).Value = 0.1454502834669897
$ws.Cells.Item(5,16).Value = 0.1454502834669897
$ws.Cells.Item(5,17).Value = 138.0729000939457
$ws.Cells.Item(5,18).Value = 1242.656100845511
$ws.Cells.Item(5,19).Value = 0.119053735944881
$ws.Cells.Item(5,20).Value = 0.119053735944881

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Fgf2"
$ws.Cells.Item(6,3).Value = "Gpc4"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 19.74619233333334
$ws.Cells.Item(6,8).Value = 59.23857700000001
$ws.Cells.Item(6,9).Value = 0.8185184181638298
$ws.Cells.Item(6,10).Value = 0.8185184181638298
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 30.83466
$ws.Cells.Item(6,14).Value = 92.50398
$ws.Cells.Item(6,15).Value = 0.641399551541635
$ws.Cells.Item(6,16).Value = 0.641399551541635
$ws.Cells.Item(6,17).Value = 608.8671268929401
$ws.Cells.Item(6,18).Value = 5479.804142036461
$ws.Cells.Item(6,19).Value = 0.5249973463388489
$ws.Cells.Item(6,20).Value = 0.5249973463388489

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Fgf2"
$ws.Cells.Item(7,3).Value = "Gpc4"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 19.74619233333334
$ws.Cells.Item(7,8).Value = 59.23857700000001
$ws.Cells.Item(7,9).Value = 0.8185184181638298
$ws.Cells.Item(7,10).Value = 0.8185184181638298
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 10.246987
$ws.Cells.Item(7,14).Value = 30.740961
$ws.Cells.Item(7,15).Value = 0.2131501649913754
$ws.Cells.Item(7,16).Value = 0.2131501649913754
$ws.Cells.Item(7,17).Value = 202.3389761391663
$ws.Cells.Item(7,18).Value = 1821.050785252497
$ws.Cells.Item(7,19).Value = 0.1744673358800999
$ws.Cells.Item(7,20).Value = 0.1744673358800999

# Row 8
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Fgf2"
$ws.Cells.Item(8,3).Value = "Gpc4"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 3.630999
$ws.Cells.Item(8,8).Value = 10.892997
$ws.Cells.Item(8,9).Value = 0.1505120332904577
$ws.Cells.Item(8,10).Value = 0.1505120332904577
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 0.6666666666666666
$ws.Cells.Item(8,13).Value = 6.992380999999999
$ws.Cells.Item(8,14).Value = 20.977143
$ws.Cells.Item(8,15).Value = 0.1454502834669897
$ws.Cells.Item(8,16).Value = 0.1454502834669897
$ws.Cells.Item(8,17).Value = 25.38932841861899
$ws.Cells.Item(8,18).Value = 228.503955767571
$ws.Cells.Item(8,19).Value = 0.02189201790729006
$ws.Cells.Item(8,20).Value = 0.02189201790729006

# Row 9
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Fgf2"
$ws.Cells.Item(9,3).Value = "Gpc4"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 3.630999
$ws.Cells.Item(9,8).Value = 10.892997
$ws.Cells.Item(9,9).Value = 0.1505120332904577
$ws.Cells.Item(9,10).Value = 0.1505120332904577
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 30.83466
$ws.Cells.Item(9,14).Value = 92.50398
$ws.Cells.Item(9,15).Value = 0.641399551541635
$ws.Cells.Item(9,16).Value = 0.641399551541635
$ws.Cells.Item(9,17).Value = 111.96061962534
$ws.Cells.Item(9,18).Value = 1007.64557662806
$ws.Cells.Item(9,19).Value = 0.09653835065411918
$ws.Cells.Item(9,20).Value = 0.0965383506541192

# Row 10
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Fgf2"
$ws.Cells.Item(10,3).Value = "Gpc4"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 3.630999
$ws.Cells.Item(10,8).Value = 10.892997
$ws.Cells.Item(10,9).Value = 0.1505120332904577
$ws.Cells.Item(10,10).Value = 0.1505120332904577
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 10.246987
$ws.Cells.Item(10,14).Value = 30.740961
$ws.Cells.Item(10,15).Value = 0.2131501649913754
$ws.Cells.Item(10,16).Value = 0.2131501649913754
$ws.Cells.Item(10,17).Value = 37.206799550013
$ws.Cells.Item(10,18).Value = 334.861195950117
$ws.Cells.Item(10,19).Value = 0.03208166472904844
$ws.Cells.Item(10,20).Value = 0.03208166472904844
